# Precision-Recall workbook update:
#  - Preserve the previous "Data" values into a new "Old Data" sheet
#    (inserted between "Data" and "Precision-Recall Graph")
#  - Update "Data" with the new values (Stemming / Stopwords&Stemming / Global Query Expansion columns)
#  - Add row-by-row delta formulas + conditional formatting on "Old Data"
#  - Add an "ir" named range scoped to "Old Data", mirroring the one on "Data"
#  - Freeze header row on "Data" and restore the reported active-cell selections

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")

# ---------------------------------------------------------------------------
# 1. Duplicate "Data" (with its current/old values) to create "Old Data",
#    positioned right after "Data" (i.e. before the chart sheet).
# ---------------------------------------------------------------------------
$wsData.Copy([Type]::Missing, $wsData)
$wsOld = $wb.Worksheets.Item(2)
$wsOld.Name = "Old Data"

# ---------------------------------------------------------------------------
# 2. Update "Data" sheet with the new values.
# ---------------------------------------------------------------------------
$newValues = @{
    "D6"  = 0.70097670182746696
    "D7"  = 0.69364336849413299
    "D8"  = 0.68896490685915002
    "D9"  = 0.58237497946083505
    "D10" = 0.54746383413609201
    "D11" = 0.53059166107807798
    "D12" = 0.50564303483807505
    "D13" = 0.50546781581429501

    "E3"  = 0.73391980507327703
    "E4"  = 0.73191980507327703
    "E5"  = 0.72358647173994395
    "E6"  = 0.68830869396216598
    "E7"  = 0.68230869396216598
    "E8"  = 0.67626023949265901
    "E9"  = 0.56727331416473703
    "E10" = 0.53469613578090702
    "E11" = 0.51792053665236704
    "E12" = 0.49169667788577898
    "E13" = 0.49169667788577898

    "I3"  = 0.67689392373980894
    "I4"  = 0.675075741921627
    "I5"  = 0.667075741921628
    "I6"  = 0.64585351969940497
    "I7"  = 0.620569969915855
    "I8"  = 0.60983187467776001
    "I9"  = 0.53920118172723897
    "I10" = 0.51887298378139401
    "I11" = 0.490916187223364
    "I12" = 0.46636176594147899
    "I13" = 0.465121999859608
}

foreach ($addr in $newValues.Keys) {
    $wsData.Range($addr).Value = $newValues[$addr]
}

# ---------------------------------------------------------------------------
# 3. Build the comparison block on "Old Data" (rows 15-24): one row per
#    original data row (3-12), each column showing Data!x - 'Old Data'!x.
# ---------------------------------------------------------------------------
$cols = @("B", "C", "D", "E", "F", "G", "H", "I")

# First set up the target number format / alignment on a single cell, then
# propagate it via copy/paste-format so the whole block shares one style.
$formatCell = $wsOld.Range("B15")
$formatCell.NumberFormat = "0.000000000000000"
$formatCell.HorizontalAlignment = -4108
$formatCell.Copy()
$deltaRange = $wsOld.Range("B15:I24")
$deltaRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($r = 15; $r -le 24; $r++) {
    $dataRow = $r - 12
    foreach ($col in $cols) {
        $wsOld.Range($col + $r).Formula = "=Data!" + $col + $dataRow + "-'Old Data'!" + $col + $dataRow
    }
}

# ---------------------------------------------------------------------------
# 4. Conditional formatting over the delta block: green when > 0, red when < 0.
# ---------------------------------------------------------------------------
$fcGreater = $deltaRange.FormatConditions.Add(1, 5, "0")
$fcGreater.Font.Color = 24832
$fcGreater.Interior.Color = 13561798

$fcLess = $deltaRange.FormatConditions.Add(1, 6, "0")
$fcLess.Font.Color = 393372
$fcLess.Interior.Color = 13551615

$fcGreater.Priority = 2
$fcLess.Priority = 1

# ---------------------------------------------------------------------------
# 5. Named range "ir" scoped to "Old Data" (mirrors the one defined on "Data").
# ---------------------------------------------------------------------------
$wsOld.Names.Add("ir", "='Old Data'!`$A`$3:`$B`$13")

# ---------------------------------------------------------------------------
# 6. Sheet view tweaks: freeze header row + active-cell selection on "Data",
#    and the reported active cell on "Old Data".
# ---------------------------------------------------------------------------
$wsData.Activate()
$wsData.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsData.Range("C3").Select()

$wsOld.Activate()
$wsOld.Range("E3").Select()

$wsData.Activate()
